$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("公告明细")
$ws2 = $wb.Worksheets.Item("公告汇总")

# Sheet1 (公告明细): rows 2..29 - full refresh of data region
# Columns: A=name, B=stock code (forced text via leading quote), C=title, D=date, E=publish time, F=attachment url
$sheet1Data = @(
    @("金丹科技", "'300829", "金丹科技:关于全资子公司取得营业执照的公告", "2023-05-17 00:00:00", "2023-05-17 11:45:48:000", "https://pdf.dfcfw.com/pdf/H2_AN202305171586712182_1.pdf?1684324097000.pdf"),
    @("金埔园林", "'301098", "金埔园林:2022年年度权益分派实施公告", "2023-05-18 00:00:00", "2023-05-17 17:16:00:000", "https://pdf.dfcfw.com/pdf/H2_AN202305171586723155_1.pdf?1684343769000.pdf"),
    @("金埔园林", "'301098", "金埔园林:江苏世纪同仁律师事务所关于金埔园林股份有限公司2022年年度股东大会的法律意见书", "2023-05-16 00:00:00", "2023-05-16 19:08:04:000", "https://pdf.dfcfw.com/pdf/H2_AN202305161586690253_1.pdf?1684264092000.pdf"),
    @("金埔园林", "'301098", "金埔园林:2022年年度股东大会决议公告", "2023-05-16 00:00:00", "2023-05-16 19:07:44:000", "https://pdf.dfcfw.com/pdf/H2_AN202305161586690224_1.pdf?1684267535000.pdf"),
    @("金埔园林", "'301098", "金埔园林:关于对外担保的进展公告", "2023-05-16 00:00:00", "2023-05-16 19:07:44:000", "https://pdf.dfcfw.com/pdf/H2_AN202305161586690223_1.pdf?1684264084000.pdf"),
    @("金埔园林", "'301098", "金埔园林:关于与兰坪白族普米族自治县人民政府签订战略合作框架协议的公告", "2023-05-16 00:00:00", "2023-05-16 19:05:59:000", "https://pdf.dfcfw.com/pdf/H2_AN202305161586689989_1.pdf?1684263959000.pdf"),
    @("九典制药", "'300705", "九典制药:300705九典制药调研活动信息20230516", "2023-05-16 00:00:00", "2023-05-16 17:24:32:000", "https://pdf.dfcfw.com/pdf/H2_AN202305161586686122_1.pdf?1684257884000.pdf"),
    @("九典制药", "'300705", "九典制药:独立董事关于第三届董事会第十八次会议相关事项的独立意见", "2023-05-16 00:00:00", "2023-05-15 17:40:58:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586642808_1.pdf?1684172462000.pdf"),
    @("九典制药", "'300705", "九典制药:监事会关于2022年限制性股票激励计划第一个归属期激励对象名单的核查意见", "2023-05-16 00:00:00", "2023-05-15 17:40:54:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586642800_1.pdf?1684225278000.pdf"),
    @("九典制药", "'300705", "九典制药:第三届监事会第十七次会议决议公告", "2023-05-16 00:00:00", "2023-05-15 17:40:54:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586642766_1.pdf?1684225564000.pdf"),
    @("九典制药", "'300705", "九典制药:关于调整2022年限制性股票激励计划相关事项的公告", "2023-05-16 00:00:00", "2023-05-15 17:40:54:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586642799_1.pdf?1684172708000.pdf"),
    @("九典制药", "'300705", "九典制药:湖南启元律师事务所关于公司2022年限制性股票激励计划相关事项调整、首次授予部分第一个归属期条件成就及部分已授予尚未归属的限制性股票作废事项的法律意见书", "2023-05-16 00:00:00", "2023-05-15 17:40:54:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586642777_1.pdf?1684172462000.pdf"),
    @("九典制药", "'300705", "九典制药:关于作废部分限制性股票的公告", "2023-05-16 00:00:00", "2023-05-15 17:40:54:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586642765_1.pdf?1684172801000.pdf"),
    @("九典制药", "'300705", "九典制药:关于2022年限制性股票激励计划第一个归属期归属条件成就的公告", "2023-05-16 00:00:00", "2023-05-15 17:40:54:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586642780_1.pdf?1684172760000.pdf"),
    @("九典制药", "'300705", "九典制药:关于延长向不特定对象发行可转换公司债券股东大会决议有效期及授权有效期的公告", "2023-05-16 00:00:00", "2023-05-15 17:40:54:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586642767_1.pdf?1684174477000.pdf"),
    @("九典制药", "'300705", "九典制药:第三届董事会第十八次会议决议公告", "2023-05-16 00:00:00", "2023-05-15 17:40:54:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586642758_1.pdf?1684235003000.pdf"),
    @("九典制药", "'300705", "九典制药:上海荣正企业咨询服务(集团)股份有限公司关于公司2022年限制性股票激励计划第一个归属期归属条件成就相关事项之独立财务顾问报告", "2023-05-16 00:00:00", "2023-05-15 17:40:54:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586642761_1.pdf?1684172462000.pdf"),
    @("九典制药", "'300705", "九典制药:关于召开2023年第一次临时股东大会的通知", "2023-05-16 00:00:00", "2023-05-15 17:40:54:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586642764_1.pdf?1684172462000.pdf"),
    @("佳禾智能", "'300793", "佳禾智能:2022年度股东大会法律意见书", "2023-05-18 00:00:00", "2023-05-18 18:34:00:000", "https://pdf.dfcfw.com/pdf/H2_AN202305181586762915_1.pdf?1684434840000.pdf"),
    @("佳禾智能", "'300793", "佳禾智能:2022年度股东大会决议公告", "2023-05-18 00:00:00", "2023-05-18 18:34:00:000", "https://pdf.dfcfw.com/pdf/H2_AN202305181586762911_1.pdf?1684434840000.pdf"),
    @("佳禾智能", "'300793", "佳禾智能:300793佳禾智能调研活动信息20230516", "2023-05-16 00:00:00", "2023-05-16 19:10:39:000", "https://pdf.dfcfw.com/pdf/H2_AN202305161586690280_1.pdf?1684264256000.pdf"),
    @("金牌厨柜", "'603180", "金牌厨柜:金牌厨柜关于2022年度业绩说明会召开情况的公告", "2023-05-17 00:00:00", "2023-05-16 17:22:52:000", "https://pdf.dfcfw.com/pdf/H2_AN202305161586686132_1.pdf?1684257887000.pdf"),
    @("金牌厨柜", "'603180", "金牌厨柜:金牌厨柜关于签订募集资金专户存储四方监管协议公告", "2023-05-17 00:00:00", "2023-05-16 17:22:52:000", "https://pdf.dfcfw.com/pdf/H2_AN202305161586686141_1.pdf?1684310250000.pdf"),
    @("金牌厨柜", "'603180", "金牌厨柜:金牌厨柜家居科技股份有限公司2022年年度股东大会会议资料", "2023-05-16 00:00:00", "2023-05-15 17:07:41:000", "https://pdf.dfcfw.com/pdf/H2_AN202305151586641321_1.pdf?1684170481000.pdf"),
    @("星帅尔", "'002860", "星帅尔:关于2022年限制性股票激励计划预留部分限制性股票登记完成的公告", "2023-05-19 00:00:00", "2023-05-18 19:54:05:000", "https://pdf.dfcfw.com/pdf/H2_AN202305181586835614_1.pdf?1684439649000.pdf"),
    @("星帅尔", "'002860", "星帅尔:关于公司取得专利证书的公告", "2023-05-19 00:00:00", "2023-05-18 18:02:55:000", "https://pdf.dfcfw.com/pdf/H2_AN202305181586760878_1.pdf?1684432991000.pdf"),
    @("星帅尔", "'002860", "星帅尔:锦天城关于星帅尔2022年年度股东大会的法律意见书", "2023-05-19 00:00:00", "2023-05-18 18:02:26:000", "https://pdf.dfcfw.com/pdf/H2_AN202305181586760849_1.pdf?1684432979000.pdf"),
    @("星帅尔", "'002860", "星帅尔:2022年年度股东大会决议公告", "2023-05-19 00:00:00", "2023-05-18 18:00:02:000", "https://pdf.dfcfw.com/pdf/H2_AN202305181586758798_1.pdf?1684432810000.pdf")
)

for ($i = 0; $i -lt $sheet1Data.Length; $i++) {
    $r = $i + 2
    $row = $sheet1Data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws1.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

# Sheet2 (公告汇总): rows 2..7 - full refresh of data region
$sheet2Data = @(
    @("金丹科技", "'300829", 1),
    @("金埔园林", "'301098", 5),
    @("九典制药", "'300705", 12),
    @("佳禾智能", "'300793", 3),
    @("金牌厨柜", "'603180", 3),
    @("星帅尔", "'002860", 4)
)

for ($i = 0; $i -lt $sheet2Data.Length; $i++) {
    $r = $i + 2
    $row = $sheet2Data[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
}
